$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(189).Insert()

$ws.Range("A189").Value = 6
$ws.Range("B189").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C189").Value = "Metropolitana"
$ws.Range("D189").Value = 44841
$ws.Range("E189").Value = 13
$ws.Range("F189").Value = 100112029
$ws.Range("G189").Value = "Orégano"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 37
$ws.Range("K189").Value = 19000
$ws.Range("L189").Value = 20000
$ws.Range("M189").Value = 19405
$ws.Range("N189").Value = "$/docena de atados"
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 6468
$ws.Range("Q189").Value = 3
$ws.Range("R189").Value = "Hortaliza"
